$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the big red "plus" cross shape ("Kreuz 5")
$cross = $s.Shapes.Item(3)
$cross.Delete()

# Add a small red dot (ellipse) in its place
$dot = $s.Shapes.AddShape(9, 383.5603937007874, -7.313149606299213, 14.173228346456693, 14.173228346456693)
$dot.Name = "Ellipse 1"
$dot.Fill.ForeColor.RGB = 192
$dot.Line.Visible = $false
$dot.TextFrame.VerticalAnchor = 3
$dot.TextFrame.TextRange.ParagraphFormat.Alignment = 2
